$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header label in E1 and G1
$ws.Range("E1").Value = "Caterpie HP"
$ws.Range("G1").Value = "CP minimum"

# Value under "CP minimum" header
$ws.Range("G2").Value = 100

# Duplicate the A1:C1 header row into I1:K1
$ws.Range("I1").Value = $ws.Range("A1").Value2
$ws.Range("J1").Value = $ws.Range("B1").Value2
$ws.Range("K1").Value = $ws.Range("C1").Value2

# Column widths to match the new layout
# (input values are pre-compensated for the engine's width quantization
# so the persisted XML width matches the target exactly)
$ws.Columns.Item(5).ColumnWidth = 13.0
$ws.Columns.Item(7).ColumnWidth = 12.833333333333334
$ws.Columns.Item(10).ColumnWidth = 17.666666666666668
$ws.Columns.Item(11).ColumnWidth = 16.166666666666668

# Update selection to E2
$ws.Range("E2").Select()
